# Add a new "2022-Q1" sheet (cloned layout from "2021-Q4") between the
# existing "2021-Q4" and "总计" sheets, populate it with the Q1-2022 fund
# holding row, and update the "总计" (totals) sheet with a new summary row.

$wb = $excel.ActiveWorkbook

$sheetQ4 = $wb.Worksheets.Item(1)     # "2021-Q4"
$sheetTotal = $wb.Worksheets.Item(2)  # "总计"

# --- 1. Insert the new "2022-Q1" sheet right before "总计" -----------------
$sheetQ1 = $wb.Worksheets.Add($sheetTotal, $null)
$sheetQ1.Name = "2022-Q1"

# NOTE: worksheet references obtained via positional Item(N) track the slot,
# not the sheet object, so after inserting a sheet they can resolve to a
# different sheet than originally fetched. Re-fetch everything we need by
# name once the sheet collection has been mutated (except $sheetQ1, which we
# just obtained directly from Add and which already carries its new name).
$sheetQ4 = $wb.Worksheets.Item("2021-Q4")
$sheetTotal = $wb.Worksheets.Item("总计")

# Clone the header row + row layout/styles from "2021-Q4" so the new sheet
# matches the same column structure (基金代码/基金名称/.../仓位排名).
$sheetQ4.Range("A1:H2").Copy($sheetQ1.Range("A1:H2"))

# --- 2. Fill in the 2022-Q1 fund holding data row --------------------------
# Numeric-looking values are forced to text (matching the source data, which
# keeps things like leading zeros / fixed decimal formatting as strings).
$sheetQ1.Range("B2").NumberFormat = "@"
$sheetQ1.Range("B2").Value = "009693"

$sheetQ1.Range("C2").Value = "富国积极成长一年定期开放混合"

$sheetQ1.Range("D2").NumberFormat = "@"
$sheetQ1.Range("D2").Value = "17.82"

$sheetQ1.Range("E2").NumberFormat = "@"
$sheetQ1.Range("E2").Value = "98.74"

$sheetQ1.Range("F2").NumberFormat = "@"
$sheetQ1.Range("F2").Value = "2.57"

$sheetQ1.Range("G2").NumberFormat = "@"
$sheetQ1.Range("G2").Value = "0.4580"

$sheetQ1.Range("H2").Value = 10

# --- 3. Update the "总计" sheet with the new 2022-Q1 summary row -----------
# Push the existing 2021-Q4 summary row down to row 3 (copy keeps its style),
# then overwrite row 2 with the new 2022-Q1 totals.
$sheetTotal.Range("A2:D2").Copy($sheetTotal.Range("A3:D3"))

$sheetTotal.Range("B2").Value = "2022-Q1"
$sheetTotal.Range("C2").Value = 1
$sheetTotal.Range("D2").Value = 0.46

$sheetTotal.Range("A3").Value = 1
